# Scheduled market-price refresh: update a couple of repriced rows on ALC,
# and clear stale price/profit figures (columns H:N) for items that no
# longer have current market data on ARM and WVR.

$wb = $excel.ActiveWorkbook

# --- ALC sheet: updated market prices for rows 87 and 90 ---
$wsAlc = $wb.Worksheets.Item("ALC")

$wsAlc.Range("H87").Value = 92475
$wsAlc.Range("J87").Value = 94950
$wsAlc.Range("L87").Value = 94950
$wsAlc.Range("N87").Value = -97446

$wsAlc.Range("H90").Value = 92475
$wsAlc.Range("J90").Value = 94950
$wsAlc.Range("L90").Value = 284850
$wsAlc.Range("N90").Value = -297330

# --- ARM sheet: clear stale price/profit columns (H:N) for rows 121-141,
#     except row 136 which still has current data ---
$wsArm = $wb.Worksheets.Item("ARM")
$wsArm.Range("H121:N135").ClearContents()
$wsArm.Range("H137:N141").ClearContents()

# --- WVR sheet: clear stale price/profit columns (H:N) for rows 119-141,
#     except row 134 which still has current data ---
$wsWvr = $wb.Worksheets.Item("WVR")
$wsWvr.Range("H119:N133").ClearContents()
$wsWvr.Range("H135:N141").ClearContents()
